$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 store a per-row "Date" label. It was computed from the wrong
# day (off by one), because of how NBA stats were shown at the time
# (4-27-2011-12 -> true ISO date 2012-04-27). Fix the label in every row
# without disturbing the cell's plain-text type/formatting.
#
# A direct $cell.Value = "2012-04-27" would be auto-recognized as a date
# literal (it looks like an ISO date) and silently reformatted as a date
# serial number with a date NumberFormat. To keep it as a literal text
# string (matching how the value was originally stored), stage the text
# in a scratch cell via a formula (whose result is always plain text),
# then Copy / PasteSpecial only the values into the real cell - this is
# the normal "convert formula to a static value" pattern and keeps the
# destination's existing General formatting/style untouched.
$scratch = $ws.Range("BZ1")

for ($row = 2; $row -le 31; $row++) {
    $scratch.Formula = '="2012-04-27"'
    $scratch.Copy()
    $ws.Range("BF$row").PasteSpecial(-4163)
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
